# Slide 4, shape "Content Placeholder 1": the bullet
#   "          Open Stego (steganography tool)"
# is split across three runs:
#   1) "          Open "   2) "Stego"   3) " (steganography tool)"
# Rename the tool to "OpenStego" by moving "Open" out of run 1 (leaving
# just the leading spaces) and prefixing it onto run 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(5, 1)

$run1 = $para.Runs(1, 1)
$run2 = $para.Runs(2, 1)

$run1.Text = "          "
$run2.Text = "OpenStego"
